# Journal de travail - ajout d'une nouvelle tache (nouvelle ligne au bas du tableau)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Je vais commencer par suivre ce tuto : https://www.youtube.com/watch?v=yWOAkwM3B9Y que j'ai trouvé sur youtube."

# Reproduit la mise en forme de la derniere ligne existante (date + texte avec
# retour a la ligne automatique) sur la nouvelle ligne 44.
$ws.Range("A43:B43").Copy()
$ws.Range("A44:B44").PasteSpecial(-4122)

$ws.Range("A44").Value = 43172
$ws.Range("B44").Value = $newText

# Meme hauteur de ligne que les autres entrees sur deux lignes du journal.
$ws.Rows.Item(44).RowHeight = 30

# La cellule active se deplace sur B45, comme apres une saisie sur la ligne 44.
$ws.Range("B45").Select() | Out-Null
